$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $val) {
    $r = $ws.Range($rangeAddr)
    $origStyle = $r.Style
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = $origStyle
}

$ws.Range("D2").Value = "66.235.22"
$ws.Range("E2").Value = "  +6.65%  "
$ws.Range("D3").Value = "3.013.39"
$ws.Range("E3").Value = "  +3.62%  "
$ws.Range("E4").Value = "  +0.05%  "
Set-TextValue "D5" "584.09"
$ws.Range("E5").Value = "  +2.92%  "
Set-TextValue "D6" "163.31"
$ws.Range("E6").Value = "  +13.64%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  +3.76%  "
$ws.Range("D9").Value = "3.009.26"
$ws.Range("E9").Value = "  +3.53%  "
Set-TextValue "D10" "6.69"
$ws.Range("E10").Value = "  -4.36%  "
$ws.Range("E11").Value = "  +4.62%  "
Set-TextValue "D12" "0.456"
$ws.Range("E12").Value = "  +5.79%  "
$ws.Range("E13").Value = "  +7.80%  "
Set-TextValue "D14" "34.71"
$ws.Range("E14").Value = "  +6.84%  "
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("D16").Value = "66.237.15"
$ws.Range("E16").Value = "  +6.75%  "
$ws.Range("D17").Value = "3.516.21"
$ws.Range("E17").Value = "  +3.69%  "
$ws.Range("E18").Value = "  +5.91%  "
$ws.Range("D19").Value = "3.014.92"
$ws.Range("E19").Value = "  +3.74%  "
Set-TextValue "D20" "456.77"
$ws.Range("E20").Value = "  +6.01%  "
Set-TextValue "D21" "13.95"
$ws.Range("E21").Value = "  +6.91%  "
$ws.Range("E22").Value = "  +5.73%  "
Set-TextValue "D23" "7.38"
$ws.Range("E23").Value = "  +7.76%  "
Set-TextValue "D24" "82.50"
$ws.Range("E24").Value = "  +4.81%  "
$ws.Range("E25").Value = "  +15.42%  "
Set-TextValue "D26" "12.35"
$ws.Range("E26").Value = "  +3.26%  "
$ws.Range("E27").Value = "  +5.18%  "
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("E29").Value = "  +17.37%  "
$ws.Range("E30").Value = "  +18.97%  "
$ws.Range("E31").Value = "  -6.10%  "
$ws.Range("E32").Value = "  +4.52%  "
Set-TextValue "D33" "27.36"
$ws.Range("E33").Value = "  +6.63%  "
$ws.Range("E34").Value = "  +5.43%  "
$ws.Range("E35").Value = "  +0.02%  "
Set-TextValue "D36" "0.995"
$ws.Range("E36").Value = "  +3.82%  "
$ws.Range("E37").Value = "  +16.59%  "
$ws.Range("E38").Value = "  +8.09%  "
$ws.Range("E39").Value = "  +3.38%  "
Set-TextValue "D40" "50.02"
$ws.Range("E40").Value = "  +2.33%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D41" "0.124"
$ws.Range("E41").Value = "  +8.97%  "
$ws.Range("B42").Value = "TheGraph"
$ws.Range("C42").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue "D42" "0.308"
$ws.Range("E42").Value = "  +16.41%  "
Set-TextValue "D43" "43.81"
$ws.Range("E43").Value = "  +7.52%  "
Set-TextValue "D44" "8.47"
$ws.Range("E44").Value = "  +4.00%  "
Set-TextValue "D45" "397.71"
$ws.Range("E45").Value = "  +15.40%  "
Set-TextValue "D46" "0.0361"
$ws.Range("E46").Value = "  +7.26%  "
$ws.Range("D47").Value = "2.801.20"
$ws.Range("E47").Value = "  +2.83%  "
Set-TextValue "D48" "134.34"
$ws.Range("E48").Value = "  +0.83%  "
Set-TextValue "D50" "23.92"
$ws.Range("E50").Value = "  +11.74%  "
